$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 63, shifting existing rows 63:96 down to 64:97.
$ws.Rows("63").Insert()

# Fill in the new row 63 with the new data record.
$ws.Range("A63").Value = 5
$ws.Range("B63").Value = "Macroferia Regional de Talca"
$ws.Range("C63").Value = "Maule"
$ws.Range("D63").Value = 44596
$ws.Range("E63").Value = 7
$ws.Range("F63").Value = 100112001
$ws.Range("G63").Value = "Berenjena"
$ws.Range("H63").Value = "Sin especificar"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 150
$ws.Range("K63").Value = 7000
$ws.Range("L63").Value = 7000
$ws.Range("M63").Value = 7000
$ws.Range("N63").Value = "$/caja 50 unidades"
$ws.Range("O63").Value = "Región del Maule"
$ws.Range("P63").Value = 140
$ws.Range("Q63").Value = 50
$ws.Range("R63").Value = "Hortaliza"
